$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 28. This pushes the old
# rows 28..144 down to 29..145 (and therefore extends the used range to
# row 145, matching the updated <dimension ref="A1:R145"/>).
$ws.Rows.Item(28).Insert()

# Populate the freshly inserted row 28 with the new record.
$ws.Range("A28").Value = 2
$ws.Range("B28").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C28").Value = "Coquimbo"
$ws.Range("D28").Value = 45274
$ws.Range("E28").Value = 4
$ws.Range("F28").Value = 100112030
$ws.Range("G28").Value = "Poroto granado"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 34000
$ws.Range("L28").Value = 36000
$ws.Range("M28").Value = 35000
$ws.Range("N28").Value = "$/malla 25 kilos"
$ws.Range("O28").Value = "Provincia de Limarí"
$ws.Range("P28").Value = 1400
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
